$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 298.55884
$ws.Range("I33").Value = 284.125
$ws.Range("J33").Value = 333.2
$ws.Range("K33").Value = 284.125
$ws.Range("L33").Value = 333.2
$ws.Range("M33").Value = -55.125
$ws.Range("N33").Value = -791.2
$ws.Range("H88").Value = 25376
$ws.Range("I88").Value = 50000
$ws.Range("J88").Value = 17168
$ws.Range("K88").Value = 50000
$ws.Range("L88").Value = 17168
$ws.Range("M88").Value = -49594
$ws.Range("N88").Value = -17980
$ws.Range("H91").Value = 25376
$ws.Range("I91").Value = 50000
$ws.Range("J91").Value = 17168
$ws.Range("K91").Value = 50000
$ws.Range("L91").Value = 17168
$ws.Range("M91").Value = -48596
$ws.Range("N91").Value = -19976
$ws.Range("H100").Value = 2898
$ws.Range("I100").Value = 2372.5
$ws.Range("K100").Value = 2372.5
$ws.Range("M100").Value = -1831.5
$ws.Range("H129").Value = 847.44446
$ws.Range("J129").Value = 849.46155
$ws.Range("L129").Value = 2548.38465
$ws.Range("N129").Value = -12548.38465
$ws.Range("H132").Value = 3460.2083
$ws.Range("I132").Value = 4110.278
$ws.Range("J132").Value = 1510
$ws.Range("K132").Value = 12330.834
$ws.Range("L132").Value = 4530
$ws.Range("M132").Value = -9800.834000000001
$ws.Range("N132").Value = -9590
$ws.Range("H137").Value = 35433.332
$ws.Range("I137").Value = 2252.2942
$ws.Range("J137").Value = 78823.92
$ws.Range("K137").Value = 6756.882599999999
$ws.Range("L137").Value = 236471.76
$ws.Range("M137").Value = -4206.882599999999
$ws.Range("N137").Value = -241571.76
$ws.Range("H138").Value = 1366.5555
$ws.Range("I138").Value = 539.25
$ws.Range("J138").Value = 3021.1667
$ws.Range("K138").Value = 1617.75
$ws.Range("L138").Value = 9063.500100000001
$ws.Range("M138").Value = 3522.25
$ws.Range("N138").Value = -19343.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23935.152
$ws.Range("I32").Value = 25287.326
$ws.Range("J32").Value = 4554
$ws.Range("K32").Value = 25287.326
$ws.Range("L32").Value = 4554
$ws.Range("M32").Value = -25000.326
$ws.Range("N32").Value = -5128
$ws.Range("H45").Value = 3112.762
$ws.Range("I45").Value = 3019
$ws.Range("J45").Value = 3176.52
$ws.Range("K45").Value = 3019
$ws.Range("L45").Value = 3176.52
$ws.Range("M45").Value = -2642
$ws.Range("N45").Value = -3930.52
$ws.Range("H74").Value = 1826.4
$ws.Range("I74").Value = 1811.8182
$ws.Range("J74").Value = 1933.3334
$ws.Range("K74").Value = 1811.8182
$ws.Range("L74").Value = 1933.3334
$ws.Range("M74").Value = -937.8181999999999
$ws.Range("N74").Value = -3681.3334
$ws.Range("H77").Value = 1826.4
$ws.Range("I77").Value = 1811.8182
$ws.Range("J77").Value = 1933.3334
$ws.Range("K77").Value = 9059.091
$ws.Range("L77").Value = 9666.666999999999
$ws.Range("M77").Value = -4691.091
$ws.Range("N77").Value = -18402.667
$ws.Range("H97").Value = 1482.6316
$ws.Range("I97").Value = 1447.6
$ws.Range("J97").Value = 1614
$ws.Range("K97").Value = 1447.6
$ws.Range("L97").Value = 1614
$ws.Range("M97").Value = -951.5999999999999
$ws.Range("N97").Value = -2606
$ws.Range("H122").Value = 1589.5883
$ws.Range("I122").Value = 1571
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 4713
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -2263
$ws.Range("N122").Value = -9850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 36223.867
$ws.Range("I134").Value = 41566
$ws.Range("K134").Value = 124698
$ws.Range("M134").Value = -122163

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8751.159
$ws.Range("I31").Value = 9681.135
$ws.Range("K31").Value = 9681.135
$ws.Range("M31").Value = -9386.135
$ws.Range("H34").Value = 8751.159
$ws.Range("I34").Value = 9681.135
$ws.Range("K34").Value = 9681.135
$ws.Range("M34").Value = -9479.135

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4108.6665
$ws.Range("I63").Value = 1708
$ws.Range("K63").Value = 5124
$ws.Range("M63").Value = -4375
$ws.Range("H66").Value = 4108.6665
$ws.Range("I66").Value = 1708
$ws.Range("K66").Value = 15372
$ws.Range("M66").Value = -11628
$ws.Range("H76").Value = 4345.8335
$ws.Range("J76").Value = 4915
$ws.Range("L76").Value = 14745
$ws.Range("N76").Value = -15511
$ws.Range("H79").Value = 4345.8335
$ws.Range("J79").Value = 4915
$ws.Range("L79").Value = 14745
$ws.Range("N79").Value = -17397
$ws.Range("H131").Value = 737.11
$ws.Range("J131").Value = 752.51044
$ws.Range("L131").Value = 2257.53132
$ws.Range("N131").Value = -12337.53132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1802.963
$ws.Range("I97").Value = 1195.25
$ws.Range("J97").Value = 3539.2856
$ws.Range("K97").Value = 1195.25
$ws.Range("L97").Value = 3539.2856
$ws.Range("M97").Value = -699.25
$ws.Range("N97").Value = -4531.2856
$ws.Range("H102").Value = 1520.75
$ws.Range("I102").Value = 1533.8
$ws.Range("J102").Value = 1488.125
$ws.Range("K102").Value = 1533.8
$ws.Range("L102").Value = 1488.125
$ws.Range("M102").Value = 88.20000000000005
$ws.Range("N102").Value = -4732.125
$ws.Range("H132").Value = 68037.25
$ws.Range("I132").Value = 7600
$ws.Range("K132").Value = 22800
$ws.Range("M132").Value = -20270

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1589.5
$ws.Range("I93").Value = 1612.5238
$ws.Range("J93").Value = 1428.3334
$ws.Range("K93").Value = 1612.5238
$ws.Range("L93").Value = 1428.3334
$ws.Range("M93").Value = -364.5237999999999
$ws.Range("N93").Value = -3924.3334
$ws.Range("H122").Value = 1310329.5
$ws.Range("I122").Value = 2181083
$ws.Range("J122").Value = 4199.1665
$ws.Range("K122").Value = 6543249
$ws.Range("L122").Value = 12597.4995
$ws.Range("M122").Value = -6540799
$ws.Range("N122").Value = -17497.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3789298.5
$ws.Range("I107").Value = 506.5
$ws.Range("K107").Value = 1519.5
$ws.Range("M107").Value = 400.5
$ws.Range("H122").Value = 1820.8948
$ws.Range("I122").Value = 1710.2667
$ws.Range("J122").Value = 2235.75
$ws.Range("K122").Value = 5130.800099999999
$ws.Range("L122").Value = 6707.25
$ws.Range("M122").Value = -2680.800099999999
$ws.Range("N122").Value = -11607.25
